$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the "Discovered by" values in column D (rows 2-49) to short names
$renameMap = @{
    "Vic - Random Platform test"   = "Victor";
    "Joan, night shift"            = "Joan";
    "Ezekiel - through client"     = "Ezekiel";
    "Night Shift, Ezekiel"         = "Ezekiel";
    "Vic - Random Test"            = "Victor";
    "Ezekiel - Random Test"        = "Ezekiel ";
    "Patricia - Call Center Alert" = "Patricia";
    "Vic - Call Center Alert"      = "Victor";
    "Gilbert - Random Test"        = "Gilbert";
    "Vic - Observation"            = "Victor";
    "Vic & Patricia"               = "Victor & Patricia";
}

for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cur = $cell.Value2
    if ($renameMap.ContainsKey($cur)) {
        $cell.Value = $renameMap[$cur]
    }
}

# Update the header for column D ("Discovered by" -> "Discovered_by")
$ws.Range("D1").Value = "Discovered_by"

# Match the saved selection/view state (cursor parked on the renamed header cell)
$ws.Range("D1").Select()
